# Apply the committed changes to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "datafile" to "population0data"
$ws.Name = "population0data"

# Trim trailing whitespace from cell A2's value
$ws.Range("A2").Value = "TestCytel - automation_nononcology"

# Update the active cell selection to D16 (matches the edited view state)
$ws.Range("D16").Select()
